# Update WALCL FRED data workbook:
#  - Append a new weekly observation row to the "Data" sheet
#  - Refresh the "SeriesInfo" metadata sheet to match the new pull

$wb = $excel.ActiveWorkbook

$dataSheet = $wb.Worksheets.Item("Data")
$seriesInfoSheet = $wb.Worksheets.Item("SeriesInfo")

# --- Data sheet: append new observation (row 112) ---
# Carry the date-column formatting (border/bold/date number format) down
# from the previous row, then fill in the new observation's values.
$dataSheet.Cells.Item(111, 1).Copy()
$dataSheet.Cells.Item(112, 1).PasteSpecial(-4122)
$dataSheet.Cells.Item(112, 1).Value = 45245
$dataSheet.Cells.Item(112, 2).Value = 7814.991

# --- SeriesInfo sheet: refresh metadata for the new pull ---
# A leading apostrophe keeps these date-shaped strings stored as plain text
# (matching the source feed), rather than letting Excel auto-convert them
# into date serial numbers.
$seriesInfoSheet.Cells.Item(3, 2).Value = "'2023-11-21"
$seriesInfoSheet.Cells.Item(4, 2).Value = "'2023-11-21"
$seriesInfoSheet.Cells.Item(7, 2).Value = "'2023-11-15"
$seriesInfoSheet.Cells.Item(14, 2).Value = "2023-11-16 15:33:02-06"
